# Add four new localized strings to the tr.xlsx localization sheet:
#   - MigrationQuestionPrompt / its Turkish message
#   - UserExpired / its Turkish message
# and widen column B to fit the longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B from 83 chars to 116 chars (Excel's ColumnWidth setter adds
# a ~0.8333 char padding compared to the raw stored "width" attribute, so we
# compensate by subtracting it up front).
$ws.Columns.Item(2).ColumnWidth = 116 - 5/6

# New row 79: MigrationQuestionPrompt
$ws.Range("A79").Value = "MigrationQuestionPrompt"
$ws.Range("B79").Value = 'Eski uygulamadaki öğrenci verilerinizi giriş ekranındakı ''Eski Sistem Verilerini Aktar'' butonunu kullanabilirsiniz.'

# New row 80: UserExpired
$ws.Range("A80").Value = "UserExpired"
$ws.Range("B80").Value = 'Kullanıcı hesabınızın süresi dolmuştur. TestOkur''u kullanmaya devam etmek için  web sitemizden lisans yenileme işlemi yapmalısınız'

# Match existing formatting conventions used in this sheet:
#  - column A data rows use the style of A2 (body-row style)
#  - column B data rows use the style of B2 (body-row style), except row 79
#    which picks up B1's (header) style, matching the source workbook.
$ws.Range("A2").Copy()
$ws.Range("A79").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("B79").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("A80").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B80").PasteSpecial(-4122)
